# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) for rows 2-13 of the
# active sheet to reflect the new weekly data snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 },
    @{ Row = 3;  D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 },
    @{ Row = 4;  D = 44907; M = 45; N = 25000; O = 25000; P = 25000; S = 5000 },
    @{ Row = 5;  D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 },
    @{ Row = 6;  D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 },
    @{ Row = 7;  D = 44931; M = 50; N = 18000; O = 18000; P = 18000; S = 3600 },
    @{ Row = 8;  D = 44902; M = 35; N = 12000; O = 12000; P = 12000; S = 2400 },
    @{ Row = 9;  D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 },
    @{ Row = 10; D = 44914; M = 56; N = 23000; O = 23000; P = 23000; S = 4600 },
    @{ Row = 11; D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 },
    @{ Row = 12; D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 },
    @{ Row = 13; D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("M$r").Value = $u.M
    $ws.Range("N$r").Value = $u.N
    $ws.Range("O$r").Value = $u.O
    $ws.Range("P$r").Value = $u.P
    $ws.Range("S$r").Value = $u.S
}
